$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Lok Sabha"
$ws.Range("B10").Value = "Zee News"
$ws.Range("B11").Value = "Zee Business"
$ws.Range("B14").Value = "Fox Life HD"
$ws.Range("B15").Value = "Movies Ok"
$ws.Range("B17").Value = "National Geographic HD"
$ws.Range("B19").Value = "Sony Mix"
$ws.Range("B21").Value = "Sony Max 2"
$ws.Range("B24").Value = "CNN News18"
$ws.Range("B25").Value = "CNBC Awaaz"
$ws.Range("B26").Value = "FYI TV18"
$ws.Range("B30").Value = "Vh1 HD"
$ws.Range("B36").Value = "MN+ HD"

$ws.Range("B75").Value = $null

$ws.Range("B36").Select()
